$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(87)
$r = $p.Range
$r.Collapse(0)
$text = @"
A tomada de perspectiva é o acto de assumir os pontos de vista de outra pessoa de modo a podermos compreendê-los melhor, mesmo que nós próprios não defendamos esses pontos de vista ou não concordemos com eles. Darvasi (2016) explica que a tomada de perspectiva muitas vezes envolve considerar ativamente aqueles que parecem inicialmente muito diferentes (um “grupo externo”), por exemplo, incorporando o seu “estado mental, pontos de vista e motivação” (p. 3). Parte do processo de tomada de perspectiva envolve abertura.

Foi demonstrado que o processo de tomada de perspectiva ajuda a reduzir o preconceito e a melhorar as atitudes em relação às pessoas que inicialmente parecem diferentes de você, em parte porque acabam parecendo mais semelhantes e menos como um “grupo externo” (Todd & Galinsky, 2014; Darvasi, 2016 ).

Um fator-chave envolvido na tomada de perspectiva é a capacidade de se identificar com uma perspectiva específica ou com um personagem.

Darvasi (2016) conclui que o “ponto de vista” de um determinado jogo digital é importante no processo de tomada de perspectiva e formação de identidade. Por exemplo, ele explica que em jogos em primeira pessoa, o jogador incorpora o avatar, mas não o vê. O jogador pode ter menos probabilidade de se envolver na tomada de perspectiva e, em vez disso, sua identidade ficará confusa com a do avatar.

Darvasi (2016) explica que em situações em que os jogadores podem passar da primeira para a terceira pessoa, ou podem assumir uma perspectiva de terceira pessoa, eles são capazes de assumir mais prontamente a perspectiva daquele personagem, pois podem ver o personagem e podem simpatizar mais facilmente com os pontos de vista, necessidades e experiências do personagem.

Outros jogos permitem uma perspectiva de terceira pessoa (ou até mesmo uma perspectiva mais distante, como a perspectiva “vista do céu”).

As histórias na literatura ensinam empatia ao fazer com que os leitores se identifiquem indiretamente com a forma como os personagens veem e interagem com um mundo fictício. Para que ocorra a tomada de perspectiva, o mundo narrativo ficcional deve ser imersivo, atraente e convincente para transportar o leitor (Johnson, 2012).

Os jogadores podem se identificar com seus avatares na tela; no entanto, os apegos mais fortes podem ser com os personagens não jogáveis (NPCs).

Esses NPCs podem até ajudar a transportar os jogadores para mundos fictícios. Por exemplo, descobriu-se que NPCs que compartilham suas histórias contribuem para a suspensão voluntária da descrença de um jogador (Harth, 2017; Ochs, Sabouret, & Corruble, 2009).

Pode parecer surpreendente que os jogadores estabeleçam ligações com objetos e personagens virtuais não humanos; no entanto, pesquisas de Turkle (2011) e Isbister (2016) sugerem que os seres humanos podem construir esses tipos de ligações com entidades não humanas e até virtuais.

Isbister (2016) argumentou que o apego surge quando se viaja por um tempo ao lado de um ser interdependente.

Bloom (2017) explica que a emoção pode influenciar as decisões e pode afetar a forma como as pessoas pensam através de decisões éticas, e até mesmo levar as pessoas a fazerem escolhas problemáticas.

Neste artigo, analisamos especificamente a narrativa, o fluxo e a imersão (transporte); tomada de perspectiva e identidade; agência, escolha e controle, e construção de relacionamento e emoção.
"@
$r.InsertAfter("`n`n" + $text)
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
